# Automatische test-sync: 2025-07-23 22:52:50
# Appends a new log row (#27) to the "Logs" sheet and bumps the
# "Retour / Terugbetaling" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 27

$logs.Cells.Item($newRow, 1).Value = "Ik heb het verkeerde artikel ontvangen."
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #17: Ik heb het verkeerde artikel ontvangen."
$logs.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,
Bedankt voor uw bericht. Het spijt me te horen dat u het verkeerde artikel heeft ontvangen. Om dit probleem op te lossen, ontvang ik graag wat meer informatie van u. Kunt u mij het ordernummer en de naam van het artikel dat u hebt ontvangen sturen? Op die manier kunnen we de situatie beter begrijpen en een passende oplossing voor u vinden.
Alvast bedankt voor uw medewerking.
Met vriendelijke groet,
[Naam]
E-mailassistent van [Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-07-23 22:52:04"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# The multi-line "Antwoord" text makes Excel auto-expand the row height;
# re-autofit so the row falls back to the sheet's default height (15) with
# no explicit per-row height override, matching the rest of the sheet.
$logs.Rows.Item($newRow).AutoFit()

# Extend conditional formatting ranges (D/G/H/I/J) that previously stopped at
# row 26 so they cover the freshly added row 27 as well.
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
    $fc = $logs.Range("$col" + "2:" + "$col" + "26").FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($logs.Range("$col" + "2:" + "$col" + "27"))
    }
}

# Update the Dashboard summary count for "Retour / Terugbetaling" (B3: 4 -> 5).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 5
